$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. thousands-dot formatted prices, trailing-zero decimals) are preserved exactly as text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.908.30'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '3.671.15'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '599.13'
$ws.Range("E5").Value = '  +3.41%  '
$ws.Range("D6").Value = '190.89'
$ws.Range("E6").Value = '  +10.99%  '
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").Value = '1.01'
$ws.Range("E8").Value = '  +1.35%  '
$ws.Range("D9").Value = '0.708'
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("D10").Value = '58.33'
$ws.Range("E10").Value = '  +14.23%  '
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  -4.37%  '
$ws.Range("D12").Value = '0.0000276'
$ws.Range("E12").Value = '  -4.17%  '
$ws.Range("D13").Value = '10.23'
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").Value = '4.269.21'
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("D15").Value = '3.677.73'
$ws.Range("E15").Value = '  -0.30%  '
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = '19.02'
$ws.Range("E17").Value = '  -1.47%  '
$ws.Range("D18").Value = '1.12'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '67.824.23'
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '12.56'
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").Value = '401.52'
$ws.Range("E21").Value = '  -0.74%  '
$ws.Range("D22").Value = '4.46'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").Value = '88.28'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").Value = '11.42'
$ws.Range("E24").Value = '  +6.09%  '
$ws.Range("D25").Value = '2.97'
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("D26").Value = '12.59'
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").Value = '6.03'
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '3.70'
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("D29").Value = '9.36'
$ws.Range("E29").Value = '  -0.80%  '
$ws.Range("D30").Value = '31.96'
$ws.Range("E30").Value = '  -1.53%  '
$ws.Range("D31").Value = '7.63'
$ws.Range("E31").Value = '  +3.22%  '
$ws.Range("D32").Value = '45.82'
$ws.Range("E32").Value = '  +6.69%  '
$ws.Range("D33").Value = '12.40'
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("D34").Value = '66.83'
$ws.Range("E34").Value = '  +3.27%  '
$ws.Range("E35").Value = '  +2.42%  '
$ws.Range("D36").Value = '617.22'
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("D38").Value = '0.400'
$ws.Range("E38").Value = '  +1.74%  '
$ws.Range("D39").Value = '0.0₃0786'
$ws.Range("E39").Value = '  -11.08%  '
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("D42").Value = '2.92'
$ws.Range("E42").Value = '  -1.74%  '
$ws.Range("D43").Value = '0.0429'
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").Value = '2.57'
$ws.Range("E44").Value = '  -7.07%  '
$ws.Range("D45").Value = '2.849.16'
$ws.Range("E45").Value = '  +1.42%  '
$ws.Range("E46").Value = '  +2.53%  '
$ws.Range("E47").Value = '  +3.96%  '
$ws.Range("D48").Value = '8.98'
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("D49").Value = '144.72'
$ws.Range("E49").Value = '  +4.21%  '
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("D51").Value = '2.51'
$ws.Range("E51").Value = '  -12.07%  '
